$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 0.1048293742242156
$ws.Range("C7").Value = 0.6683285205126266
$ws.Range("D7").Value = 0.9396403565264032
$ws.Range("E7").Value = 0.9693504817796312
$ws.Range("F7").Value = 0.9766011666972538
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = 0.07359833454455669
$ws.Range("C8").Value = 0.6734286924824703
$ws.Range("D8").Value = 0.9899493426166667
$ws.Range("E8").Value = 0.9949619804880319
$ws.Range("F8").Value = 1.005922838153142
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = 0.134
$ws.Range("C9").Value = 0.7460000000000001
$ws.Range("D9").Value = 1.49462
$ws.Range("E9").Value = 1.222546522632165
$ws.Range("F9").Value = 1.246749035650213
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.2623076923076923
$ws.Range("C10").Value = 0.6884615384615385
$ws.Range("D10").Value = 1.410730769230769
$ws.Range("E10").Value = 1.187741878200297
$ws.Range("F10").Value = 1.205716618489836
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = -0.06799999999999998
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = 0.37436
$ws.Range("E11").Value = 0.6118496547355404
$ws.Range("F11").Value = 0.6798308613177252
$ws.Range("G11").Value = 5
